$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.494056333333333
$ws.Range("H2").Value = 4.482169
$ws.Range("I2").Value = 0.2373598341615283
$ws.Range("J2").Value = 0.2373598341615284
$ws.Range("M2").Value = 7.574702666666667
$ws.Range("N2").Value = 22.724108
$ws.Range("O2").Value = 0.2360813295275979
$ws.Range("P2").Value = 0.2360813295275979
$ws.Range("Q2").Value = 11.31703249225022
$ws.Range("R2").Value = 101.853292430252
$ws.Range("S2").Value = 0.05603622522530377
$ws.Range("T2").Value = 0.05603622522530378
$ws.Range("G3").Value = 1.494056333333333
$ws.Range("H3").Value = 4.482169
$ws.Range("I3").Value = 0.2373598341615283
$ws.Range("J3").Value = 0.2373598341615284
$ws.Range("O3").Value = 0.5879438355171306
$ws.Range("P3").Value = 0.5879438355171307
$ws.Range("Q3").Value = 28.18426812268422
$ws.Range("R3").Value = 253.658413104158
$ws.Range("S3").Value = 0.139554251294639
$ws.Range("T3").Value = 0.1395542512946391
$ws.Range("G4").Value = 1.494056333333333
$ws.Range("H4").Value = 4.482169
$ws.Range("I4").Value = 0.2373598341615283
$ws.Range("J4").Value = 0.2373598341615284
$ws.Range("O4").Value = 0.1759748349552714
$ws.Range("P4").Value = 0.1759748349552714
$ws.Range("Q4").Value = 8.435707003989776
$ws.Range("R4").Value = 75.92136303590799
$ws.Range("S4").Value = 0.04176935764158553
$ws.Range("T4").Value = 0.04176935764158554
$ws.Range("I5").Value = 0.4879832509286579
$ws.Range("J5").Value = 0.4879832509286579
$ws.Range("M5").Value = 7.574702666666667
$ws.Range("N5").Value = 22.724108
$ws.Range("O5").Value = 0.2360813295275979
$ws.Range("P5").Value = 0.2360813295275979
$ws.Range("Q5").Value = 23.26645671093334
$ws.Range("R5").Value = 209.3981103984
$ws.Range("S5").Value = 0.115203734666437
$ws.Range("T5").Value = 0.115203734666437
$ws.Range("I6").Value = 0.4879832509286579
$ws.Range("J6").Value = 0.4879832509286579
$ws.Range("O6").Value = 0.5879438355171306
$ws.Range("P6").Value = 0.5879438355171307
$ws.Range("S6").Value = 0.2869067442191135
$ws.Range("T6").Value = 0.2869067442191136
$ws.Range("I7").Value = 0.4879832509286579
$ws.Range("J7").Value = 0.4879832509286579
$ws.Range("O7").Value = 0.1759748349552714
$ws.Range("P7").Value = 0.1759748349552714
$ws.Range("S7").Value = 0.08587277204310735
$ws.Range("T7").Value = 0.08587277204310735
$ws.Range("I8").Value = 0.2746569149098138
$ws.Range("J8").Value = 0.2746569149098139
$ws.Range("M8").Value = 7.574702666666667
$ws.Range("N8").Value = 22.724108
$ws.Range("O8").Value = 0.2360813295275979
$ws.Range("P8").Value = 0.2360813295275979
$ws.Range("Q8").Value = 13.095312613592
$ws.Range("R8").Value = 117.857813522328
$ws.Range("S8").Value = 0.06484136963585718
$ws.Range("T8").Value = 0.0648413696358572
$ws.Range("I9").Value = 0.2746569149098138
$ws.Range("J9").Value = 0.2746569149098139
$ws.Range("O9").Value = 0.5879438355171306
$ws.Range("P9").Value = 0.5879438355171307
$ws.Range("S9").Value = 0.1614828400033781
$ws.Range("T9").Value = 0.1614828400033782
$ws.Range("I10").Value = 0.2746569149098138
$ws.Range("J10").Value = 0.2746569149098139
$ws.Range("O10").Value = 0.1759748349552714
$ws.Range("P10").Value = 0.1759748349552714
$ws.Range("S10").Value = 0.0483327052705785
$ws.Range("T10").Value = 0.04833270527057851
